# The commit adds one new weekly price record for
# "Bruselas (repollito)" at Vega Central Mapocho de Santiago.
# The new record is inserted as row 21, pushing the previously
# existing rows 21-59 down to rows 22-60 (dates/prices keep their
# original relative order, only shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; Excel shifts rows 21:59 down to 22:60
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A21").Value2 = 9
$ws.Range("B21").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value2 = "Metropolitana"
$ws.Range("D21").Value2 = 44791
$ws.Range("E21").Value2 = 13
$ws.Range("F21").Value2 = 100112035
$ws.Range("G21").Value2 = "Bruselas (repollito)"
$ws.Range("H21").Value2 = "Sin especificar"
$ws.Range("I21").Value2 = "Primera"
$ws.Range("J21").Value2 = 34
$ws.Range("K21").Value2 = 20000
$ws.Range("L21").Value2 = 20000
$ws.Range("M21").Value2 = 20000
$ws.Range("N21").Value2 = "`$/malla 15 kilos"
$ws.Range("O21").Value2 = "Hijuelas"
$ws.Range("P21").Value2 = 1333
$ws.Range("Q21").Value2 = 15
$ws.Range("R21").Value2 = "Hortaliza"
